$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")
$outputsWs = $wb.Worksheets.Item("Outputs")

# ---------------------------------------------------------------
# Insert 4 new rows (14-17) above the "Problem Definition" section,
# shifting everything below down by 4 rows.
# ---------------------------------------------------------------
$ws.Rows("14:17").Insert()

# Carry over the formatting used by the other value cells in this
# section (fill + right-aligned text) onto the new B13:B16 cells.
$ws.Range("B10").Copy()
$ws.Range("B13:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New rows of data under the "Cluster/Server Options" section
$ws.Range("A13").Value = "Allow Multiple Jobs"
$ws.Range("B13").Value = $true

$ws.Range("A14").Value = "Use Server As Worker"
$ws.Range("B14").Value = $true

$ws.Range("A15").Value = "Simulate Data Point Filename"
$ws.Range("B15").Value = "simulate_data_point.rb"

$ws.Range("A16").Value = "Run Data Point Filename"
$ws.Range("B16").Value = "run_openstudio_workflow.rb"

# ---------------------------------------------------------------
# Column widths for A and B (bestFit replaced with explicit widths)
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 25
$ws.Columns("B").ColumnWidth = 26

# ---------------------------------------------------------------
# Selection / active sheet changes: Setup becomes the active tab,
# with B17 selected (Outputs loses tabSelected).
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("B17").Select()

Write-Host "done"
